# Populate the "user_details" sheet (sheet2) with the Name/Surname/Postalcode
# table that was read in from the external data file, matching the order in
# which the values were typed in by the author (column-by-column: Surname,
# Postalcode, then Name, then the two data rows) so the shared-string table
# ends up in the same order as the authored workbook.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("B1").Value = "Surname"
$ws2.Range("C1").Value = "Postalcode"
$ws2.Range("A1").Value = "Name"

$ws2.Range("A2").Value = "Khulekani"
$ws2.Range("B2").Value = "Hlengwa"
$ws2.Range("C2").Value = 36845

$ws2.Range("A3").Value = "Mlungu"
$ws2.Range("B3").Value = "Wamodimo"
$ws2.Range("C3").Value = 33502

# Match the existing look of the "login" sheet: a highlighted/bordered
# header row and a bordered body, by copying the formats over instead of
# re-creating new styles.
$ws1.Range("A1:B1").Copy()
$ws2.Range("A1:C1").PasteSpecial(-4122)
$ws1.Range("A2:B3").Copy()
$ws2.Range("A2:C3").PasteSpecial(-4122)

# Auto-size the new columns to fit their contents.
$ws2.Columns.Item(1).AutoFit()
$ws2.Columns.Item(2).AutoFit()
$ws2.Columns.Item(3).AutoFit()

# user_details becomes the active/selected sheet, with C5 selected.
$ws2.Activate()
$null = $ws2.Range("C5").Select()
